$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the disclaimer text in A4 into two cells: keep the first part in A4
# (with trailing space preserved) and move the second part into a new A5 cell.
$ws.Range("A4").Value = "Links have been provided where we purchased these items. No guarantee that the prices will be the same after the date of completion for this project. "
$ws.Range("A5").Value = "We cannot take any responsibilty for purchases that you make from these URLs provided"

# Update the selected cell/range shown in the sheet view.
$ws.Range("B24").Select()
